# Commit: "Update IDEAS FOR PROBLEM STATEMENT 4.pptx"
#
# The second slide in the deck (sldId 258) only contains an empty
# title/body placeholder pair plus a small "Abhishek" text box - it was
# removed from the presentation, shifting every following slide (and its
# associated notes page) up by one position while leaving their content
# untouched.

$p = $ppt.ActivePresentation

$s = $p.Slides.Item(2)
$s.Delete()
